$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
# Row 70
$ws.Range("H70").Value = 1485
$ws.Range("I70").Value = 1537.5
$ws.Range("J70").Value = 1275
$ws.Range("K70").Value = 4612.5
$ws.Range("L70").Value = 3825
$ws.Range("M70").Value = -4342.5
$ws.Range("N70").Value = -4365
# Row 73
$ws.Range("H73").Value = 1485
$ws.Range("I73").Value = 1537.5
$ws.Range("J73").Value = 1275
$ws.Range("K73").Value = 4612.5
$ws.Range("L73").Value = 3825
$ws.Range("M73").Value = -3676.5
$ws.Range("N73").Value = -5697

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12446.5
$ws.Range("I32").Value = 12446.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 12446.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -12159.5
# Row 61
$ws.Range("H61").Value = 4133.3335
$ws.Range("I61").Value = 4200
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 4200
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3988
$ws.Range("N61").Value = -4424
# Row 88
$ws.Range("H88").Value = 1107.3334
$ws.Range("I88").Value = 999.5
$ws.Range("J88").Value = 1161.25
$ws.Range("K88").Value = 999.5
$ws.Range("L88").Value = 1161.25
$ws.Range("M88").Value = -593.5
$ws.Range("N88").Value = -1973.25
# Row 91
$ws.Range("H91").Value = 1107.3334
$ws.Range("I91").Value = 999.5
$ws.Range("J91").Value = 1161.25
$ws.Range("K91").Value = 999.5
$ws.Range("L91").Value = 1161.25
$ws.Range("M91").Value = 404.5
$ws.Range("N91").Value = -3969.25
# Row 132
$ws.Range("H132").Value = 9004.556
$ws.Range("I132").Value = 7321.8
$ws.Range("J132").Value = 11108
$ws.Range("K132").Value = 21965.4
$ws.Range("L132").Value = 33324
$ws.Range("M132").Value = -19435.4
$ws.Range("N132").Value = -38384
# Row 136
$ws.Range("H136").Value = 4133.3335
$ws.Range("I136").Value = 4200
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 12600
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -10050
$ws.Range("N136").Value = -17100

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 92
$ws.Range("H92").Value = 19500
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 19500
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 19500
$ws.Range("N92").Value = -24492
# Row 97
$ws.Range("H97").Value = 100000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 100000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 100000
$ws.Range("N97").Value = -101982
# Row 99
$ws.Range("H99").Value = 11000
$ws.Range("I99").Value = 8000
$ws.Range("J99").Value = 14000
$ws.Range("K99").Value = 8000
$ws.Range("L99").Value = 14000
$ws.Range("M99").Value = -6502
$ws.Range("N99").Value = -16996
# Row 102
$ws.Range("H102").Value = 50000
$ws.Range("I102").Value = 50000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 50000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -47566
$ws.Range("N102").ClearContents()
# Row 104
$ws.Range("H104").Value = 100000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 100000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 100000
$ws.Range("N104").Value = -105242
# Row 126
$ws.Range("H126").Value = 11000
$ws.Range("I126").Value = 8000
$ws.Range("J126").Value = 14000
$ws.Range("K126").Value = 24000
$ws.Range("L126").Value = 42000
$ws.Range("M126").Value = -21530
$ws.Range("N126").Value = -46940

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 5269.6924
$ws.Range("I117").Value = 389.2857
$ws.Range("J117").Value = 10963.5
$ws.Range("K117").Value = 1167.8571
$ws.Range("L117").Value = 32890.5
$ws.Range("M117").Value = 2274.1429
# Row 139
$ws.Range("H139").Value = 2094.8333
$ws.Range("I139").Value = 1189.6666
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 3568.9998
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 1571.0002

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 1300
$ws.Range("I6").Value = 1300
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("M6").Value = -1187
# Row 16
$ws.Range("H16").Value = 1300
$ws.Range("I16").Value = 1300
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1300
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("M16").Value = -1050
# Row 58
$ws.Range("H58").Value = 10031
$ws.Range("I58").Value = 10031
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 10031
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -9754
# Row 63
$ws.Range("H63").Value = 46000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 46000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 46000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -47372
# Row 66
$ws.Range("H66").Value = 46000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 46000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 138000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -144864

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 823
$ws.Range("I61").Value = 802.5
$ws.Range("J61").Value = 905
$ws.Range("K61").Value = 802.5
$ws.Range("L61").Value = 905
$ws.Range("M61").Value = -600.5
$ws.Range("N61").Value = -1309
# Row 113
$ws.Range("H113").Value = 823
$ws.Range("I113").Value = 802.5
$ws.Range("J113").Value = 905
$ws.Range("K113").Value = 802.5
$ws.Range("L113").Value = 905
$ws.Range("M113").Value = 1367.5
$ws.Range("N113").Value = -5245
# Row 122
$ws.Range("H122").Value = 8174.5
$ws.Range("I122").Value = 7566
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 22698
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -20248
# Row 132
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -470
$ws.Range("N132").ClearContents()
# Row 136
$ws.Range("H136").Value = 163884.33
$ws.Range("I136").Value = 5852
$ws.Range("J136").Value = 479949
$ws.Range("K136").Value = 17556
$ws.Range("L136").Value = 1439847
$ws.Range("M136").Value = -15006
$ws.Range("N136").Value = -1444947

# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 102
$ws.Range("H102").Value = 100000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 100000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 100000
$ws.Range("N102").Value = -106490
# Row 113
$ws.Range("H113").Value = 221
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 221
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 663
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5003
# Row 136
$ws.Range("H136").Value = 2228.5715
$ws.Range("I136").Value = 933.3333
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 2799.9999
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -249.9998999999998
$ws.Range("N136").Value = -35100
